$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 320; this shifts existing rows 320..420 down to 321..421
# and automatically extends the sheet dimension (A1:R420 -> A1:R421).
$ws.Rows(320).Insert()

# Populate the newly inserted row 320 with the new data record.
$ws.Range("A320").Value = 5
$ws.Range("B320").Value = "Macroferia Regional de Talca"
$ws.Range("C320").Value = "Maule"
$ws.Range("D320").Value = 44876
$ws.Range("E320").Value = 7
$ws.Range("F320").Value = 100114014
$ws.Range("G320").Value = "Betarraga"
$ws.Range("H320").Value = "Sin especificar"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 5000
$ws.Range("K320").Value = 800
$ws.Range("L320").Value = 800
$ws.Range("M320").Value = 800
$ws.Range("N320").Value = "$/paquete 5 unidades"
$ws.Range("O320").Value = "Región del Maule"
$ws.Range("P320").Value = 160
$ws.Range("Q320").Value = 5
$ws.Range("R320").Value = "Hortaliza"

# Make sure the D320 cell keeps/gets the date number format used by the rest of column D.
$ws.Range("D320").NumberFormat = "YYYY-MM-DD HH:MM:SS"
